$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update values and reasoning text ---
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 78.53
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = "`n`n`nReasoning: The candidate's projects demonstrate a strong knowledge of the skills required for the job, such as ReactJS, HTML, CSS, Flutter, Dart, Firebase, NodeJS, ExpressJS, Socket.IO, WebRTC, JS, and Docker. The projects also show a solid understanding of web development, user experience, and brand consistency. The score of 78.53 reflects the candidate's aptitude for the job."

# --- Row 3: update values and reasoning text (replaces old row3/row4 content) ---
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 76.81
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = "`n`n`nReasoning: The candidate has a strong set of relevant skills and experience for the job, including ReactJS, JavaScript, CSS, Frontend Development, NextJS, Django Rest Framework, Pytorch, Tensorflow, Keras, and Sklearn. The candidate's projects demonstrate a sound understanding of the technologies used, and show a good level of creativity and problem solving. The score of 76.81 reflects the candidate's aptitude for the job."

# --- Row 4: delete entirely (shrinks dimension to A1:E3) ---
$ws.Range("A4:E4").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)
